$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first 9 data rows (rows 2-10), shifting remaining data up.
$ws.Range("A2:C10").EntireRow.Delete() | Out-Null

# After the deletion, former rows 11-22 are now rows 2-13.
# Append new rows 14-21 with fresh data.
$newData = @(
    @(-2.109282225370405, -3.67548027634622, -7.529284000396693),
    @(-1.644850492477401, -7.059904575347879, 1.617063522338856),
    @(1.835043907165529, -2.18121553957462, -0.952013134956361),
    @(8.433930218219764, -3.52629014849663, 2.835070371627813),
    @(-5.246673464775119, -3.962655484676361, -5.154342770576502),
    @(-7.417413711547821, -2.038821458816521, -6.523755788803086),
    @(0.6842400431633051, -3.085346877574939, -2.034696012735365),
    @(1.740720510482789, -5.68582010269165, -1.830426752567291)
)

$row = 14
foreach ($rec in $newData) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row++
}
